$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D price values keep their exact text representation
# (some look numeric, e.g. "0.540", "5.30", "1.750.42") instead of being
# auto-converted to numbers by Excel, matching the source data feed.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Updated crypto price / 1h volume data ---
$ws.Range('D2').Value = '29.640.28'
$ws.Range('E2').Value = '  +3.51%  '
$ws.Range('D3').Value = '1.608.31'
$ws.Range('E3').Value = '  +2.84%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '212.66'
$ws.Range('E5').Value = '  +1.13%  '
$ws.Range('D6').Value = '0.521'
$ws.Range('E6').Value = '  +2.29%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '27.31'
$ws.Range('E8').Value = '  +10.04%  '
$ws.Range('D9').Value = '43.64'
$ws.Range('E9').Value = '  -1.40%  '
$ws.Range('E10').Value = '  +2.51%  '
$ws.Range('E11').Value = '  +2.43%  '
$ws.Range('D12').Value = '0.0909'
$ws.Range('E12').Value = '  +1.41%  '
$ws.Range('D13').Value = '1.839.34'
$ws.Range('E13').Value = '  +2.87%  '
$ws.Range('D14').Value = '1.604.00'
$ws.Range('E14').Value = '  +2.91%  '
$ws.Range('D15').Value = '29.626.42'
$ws.Range('E15').Value = '  +3.30%  '
$ws.Range('E16').Value = '  +4.22%  '
$ws.Range('E17').Value = '  +2.60%  '
$ws.Range('D18').Value = '63.53'
$ws.Range('E18').Value = '  +3.23%  '
$ws.Range('D19').Value = '241.07'
$ws.Range('E19').Value = '  +5.90%  '
$ws.Range('E20').Value = '  +4.05%  '
$ws.Range('D21').Value = '0.0₃0695'
$ws.Range('E21').Value = '  +2.05%  '
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('E23').Value = '  +2.13%  '
$ws.Range('D24').Value = '9.24'
$ws.Range('E24').Value = '  +1.83%  '
$ws.Range('D25').Value = '2.09'
$ws.Range('E25').Value = '  +0.74%  '
$ws.Range('D26').Value = '155.32'
$ws.Range('E26').Value = '  +2.41%  '
$ws.Range('D27').Value = '15.36'
$ws.Range('E27').Value = '  +3.92%  '
$ws.Range('E28').Value = '  +1.72%  '
$ws.Range('D29').Value = '6.42'
$ws.Range('E29').Value = '  +2.88%  '
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('E31').Value = '  +3.86%  '
$ws.Range('E32').Value = '  +0.92%  '
$ws.Range('D33').Value = '3.22'
$ws.Range('E33').Value = '  +1.17%  '
$ws.Range('D34').Value = '1.435.36'
$ws.Range('E34').Value = '  +2.29%  '
$ws.Range('E35').Value = '  +4.16%  '
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('E37').Value = '  +4.80%  '
$ws.Range('E38').Value = '  +5.48%  '
$ws.Range('E39').Value = '  +0.33%  '
$ws.Range('E40').Value = '  +1.84%  '
$ws.Range('D41').Value = '0.540'
$ws.Range('E41').Value = '  +4.46%  '
$ws.Range('E42').Value = '  +2.09%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').Value = '0.0490'
$ws.Range('E43').Value = '  +5.73%  '
$ws.Range('B44').Value = 'BitcoinSV'
$ws.Range('C44').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D44').Value = '53.92'
$ws.Range('E44').Value = '  +26.43%  '
$ws.Range('D45').Value = '0.802'
$ws.Range('E45').Value = '  +4.56%  '
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').Value = '66.04'
$ws.Range('E47').Value = '  +3.26%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = '0.947'
$ws.Range('E48').Value = '  +10.49%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = '5.30'
$ws.Range('E49').Value = '  +1.61%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '1.749.77'
$ws.Range('E50').Value = '  +3.05%  '
$ws.Range('D51').Value = '86.86'
$ws.Range('E51').Value = '  +2.42%  '
